$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.722.03"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "2.738.77"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'579.96"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'159.00"
$ws.Range("E6").Value = "  +10.24%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "2.760.36"
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = "  +3.90%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "3.245.12"
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("D15").Value = "'27.70"
$ws.Range("E15").Value = "  +5.35%  "
$ws.Range("D16").Value = "63.748.23"
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("E17").Value = "  +7.73%  "
$ws.Range("D18").Value = "2.752.07"
$ws.Range("E18").Value = "  +3.70%  "
$ws.Range("D19").Value = "'12.17"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("D21").Value = "'364.19"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "'7.02"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").Value = "'0.541"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'66.68"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("D27").Value = "'8.62"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +11.79%  "
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "'7.26"
$ws.Range("E31").Value = "  +6.13%  "
$ws.Range("E32").Value = "  +14.90%  "
$ws.Range("D33").Value = "'173.69"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("D36").Value = "'4.96"
$ws.Range("E36").Value = "  +6.27%  "
$ws.Range("E37").Value = "  +7.80%  "
$ws.Range("E38").Value = "  +7.14%  "
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("D41").Value = "'339.26"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +15.94%  "
$ws.Range("D43").Value = "'39.68"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "'22.56"
$ws.Range("E44").Value = "  +6.60%  "
$ws.Range("D45").Value = "'21.91"
$ws.Range("E45").Value = "  +6.49%  "
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("D47").Value = "'0.647"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("D49").Value = "'138.01"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "  -0.07%  "
